$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The diff removes ", full_name" from the "editNameUser" output line
# (paragraph 33: "output: user_id, full_name, status, and message" ->
#  "output: user_id, status, and message"). The pre-existing "_GoBack"
# bookmark (which sat in the editPassUser output paragraph before the edit)
# ends up right after "user_id," in this newly-edited paragraph - exactly
# where Word leaves it after the last edit made in the document. The two
# editPassUser paragraphs keep the same wording but some of their runs
# coalesce as a result (no visible text change there).
# ---------------------------------------------------------------------------

# --- 1) editNameUser's "output:" paragraph (#33) loses ", full_name" ------
$p33 = $d.Paragraphs.Item(33)
$p33Start = $p33.Range.Start
$t33 = $p33.Range.Text.TrimEnd([char]13)
if ($t33 -ne "output: user_id, full_name, status, and message") {
    throw "Unexpected text in paragraph 33: [$t33]"
}
$prefixLen33 = ("output: user_id").Length
$delStart33 = $p33Start + $prefixLen33
$delLen33 = (", full_name,").Length
$delEnd33 = $delStart33 + $delLen33
$delRange33 = $d.Range($delStart33, $delEnd33)
$delRange33.Text = ","

# --- 2) editPassUser's "required input:" paragraph (#35): the ", password, "
#        and "and " runs coalesce into a single ", password, and " run ----
$p35 = $d.Paragraphs.Item(35)
$p35Start = $p35.Range.Start
$t35 = $p35.Range.Text.TrimEnd([char]13)
if ($t35 -ne "required input: user_id, password, and new_pass") {
    throw "Unexpected text in paragraph 35: [$t35]"
}
$prefixLen35 = ("required input: user_id").Length
$runAStart35 = $p35Start + $prefixLen35
$runALen35 = (", password, ").Length
$runAEnd35 = $runAStart35 + $runALen35
$runBStart35 = $runAEnd35
$runBLen35 = ("and ").Length
$runBEnd35 = $runBStart35 + $runBLen35
$runBRange35 = $d.Range($runBStart35, $runBEnd35)
$runBRange35.Delete()
$runARange35 = $d.Range($runAStart35, $runAEnd35)
$runARange35.Text = ", password, and "

# --- 3) editPassUser's "output:" paragraph (#36): the "," and
#        " status, and message" runs (straddling the old bookmark)
#        coalesce into a single ", status, and message" run --------------
$p36 = $d.Paragraphs.Item(36)
$p36Start = $p36.Range.Start
$t36 = $p36.Range.Text.TrimEnd([char]13)
if ($t36 -ne "output: user_id, status, and message") {
    throw "Unexpected text in paragraph 36: [$t36]"
}
$prefixLen36 = ("output: user_id").Length
$runAStart36 = $p36Start + $prefixLen36
$runALen36 = (",").Length
$runAEnd36 = $runAStart36 + $runALen36
$runBStart36 = $runAEnd36
$runBLen36 = (" status, and message").Length
$runBEnd36 = $runBStart36 + $runBLen36
$runBRange36 = $d.Range($runBStart36, $runBEnd36)
$runBRange36.Delete()
$runARange36 = $d.Range($runAStart36, $runAEnd36)
$runARange36.Text = ", status, and message"

# --- 4) Relocate the "_GoBack" bookmark to right after "user_id," in the
#        newly-edited paragraph #33 (adding one with the same name moves it,
#        removing it from its old spot automatically) ---------------------
$bmStart = $p33Start + ("output: user_id,").Length
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Host "edit applied"
